# Insert a new product row ("Gelatina ... sabor naranja") right above the
# existing "Pan ... panchos" row on the "Artículos" sheet, shifting every
# row below it down by one (matches the target: old row 8 -> new row 9,
# ..., old row 70 -> new row 71).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# Insert a blank row at position 8 (pushes everything from old row 8 down).
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new article's data.
$ws.Cells.Item(8, 1).Value = 7790070432575
$ws.Cells.Item(8, 2).Value = "Gelatina"
$ws.Cells.Item(8, 3).Value = "en polvo"
$ws.Cells.Item(8, 4).Value = "sabor naranja"
$ws.Cells.Item(8, 5).Value = "Exquisita"
$ws.Cells.Item(8, 6).Value = 40
$ws.Cells.Item(8, 7).Value = "gr."
$ws.Cells.Item(8, 8).Value = "sobre"
$ws.Cells.Item(8, 9).Value = "Gelatinas"
$ws.Cells.Item(8, 10).Value = "Argentina"
$ws.Cells.Item(8, 11).Value = 12
$ws.Cells.Item(8, 12).Value = $false
$ws.Cells.Item(8, 13).Value = $true
# Column 14 (Imagen) intentionally left blank - no image yet for this article.
$ws.Cells.Item(8, 15).Value = $false
$ws.Cells.Item(8, 16).Value = $true
